$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.524.76"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "2.000.87"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -9.57%  "
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.82"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.371"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.13"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0978"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("D13").Value = "2.296.40"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("E16").Value = "  -5.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("D18").Value = "1.999.20"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "36.464.90"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.83"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  -8.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -5.19%  "
$ws.Range("E34").Value = "  -6.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.64"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "1.453.19"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0930"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("E44").Value = "  -4.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.15"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.77"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +23.69%  "
$ws.Range("E51").Value = "  -2.41%  "
